# Update stats for 2025-07 (row 20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 6170
$ws.Range("C20").Value = 977
$ws.Range("D20").Value = 5571514
$ws.Range("E20").Value = 903.0006482982171
$ws.Range("F20").Value = 6.581447572983246
$ws.Range("G20").Value = 3.715498938428885
$ws.Range("H20").Value = 26.03316282385057
